$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared strings with uniform run formatting) ---
$ws.Range("A8").Value = "Volume 30   Number  37"
$ws.Range("C9").Value = "Report Covering the Week  9/11/2023  Through  9/17/2023"

# --- Data table updates (rows 14-29) ---
$ws.Range("C14").Copy()
$ws.Range("F14").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("F14").PasteSpecial(-4122)
$ws.Range("N14").Value = -57.142857142857
$ws.Range("C16").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C15").Value = 1
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 50
$ws.Range("I15").Value = 24
$ws.Range("K15").Value = 9.090909090909
$ws.Range("L15").Value = -14.285714285714
$ws.Range("M15").Value = 71.428571428571
$ws.Range("N15").Value = -46.666666666666
$ws.Range("C16").Value = 8
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 166.666666666667
$ws.Range("F16").Value = 35
$ws.Range("G16").Value = 21
$ws.Range("H16").Value = 66.666666666666
$ws.Range("I16").Value = 293
$ws.Range("J16").Value = 316
$ws.Range("K16").Value = -7.278481012658
$ws.Range("L16").Value = 30.803571428571
$ws.Range("M16").Value = 12.692307692307
$ws.Range("N16").Value = -63.14465408805
$ws.Range("C17").Value = 24
$ws.Range("D17").Value = 18
$ws.Range("E17").Value = 33.333333333333
$ws.Range("F17").Value = 69
$ws.Range("G17").Value = 57
$ws.Range("H17").Value = 21.052631578947
$ws.Range("I17").Value = 508
$ws.Range("J17").Value = 460
$ws.Range("K17").Value = 10.434782608695
$ws.Range("L17").Value = 53.012048192771
$ws.Range("M17").Value = 87.453874538745
$ws.Range("N17").Value = 0
$ws.Range("C18").Value = 10
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = 42.857142857142
$ws.Range("F18").Value = 26
$ws.Range("G18").Value = 21
$ws.Range("H18").Value = 23.809523809523
$ws.Range("I18").Value = 179
$ws.Range("J18").Value = 184
$ws.Range("K18").Value = -2.717391304347
$ws.Range("L18").Value = 22.602739726027
$ws.Range("M18").Value = 5.294117647058
$ws.Range("N18").Value = -81.138040042149
$ws.Range("C19").Value = 21
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = 200
$ws.Range("F19").Value = 51
$ws.Range("G19").Value = 31
$ws.Range("H19").Value = 64.516129032258
$ws.Range("I19").Value = 369
$ws.Range("J19").Value = 323
$ws.Range("K19").Value = 14.241486068111
$ws.Range("L19").Value = 13.888888888888
$ws.Range("M19").Value = 57.021276595744
$ws.Range("N19").Value = 10.81081081081
$ws.Range("C20").Value = 9
$ws.Range("D20").Value = 8
$ws.Range("E20").Value = 12.5
$ws.Range("F20").Value = 28
$ws.Range("G20").Value = 21
$ws.Range("H20").Value = 33.333333333333
$ws.Range("I20").Value = 249
$ws.Range("J20").Value = 206
$ws.Range("K20").Value = 20.873786407767
$ws.Range("L20").Value = 91.538461538461
$ws.Range("M20").Value = 203.658536585366
$ws.Range("N20").Value = -36.479591836734
$ws.Range("C21").Value = 73
$ws.Range("D21").Value = 43
$ws.Range("E21").Value = 69.767441860465
$ws.Range("F21").Value = 212
$ws.Range("G21").Value = 153
$ws.Range("H21").Value = 38.562091503268
$ws.Range("I21").Value = 1634
$ws.Range("J21").Value = 1523
$ws.Range("K21").Value = 7.288246881155
$ws.Range("L21").Value = 36.96563285834
$ws.Range("M21").Value = 57.874396135265
$ws.Range("N21").Value = -46.426229508196
$ws.Range("C14").Copy()
$ws.Range("C22").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C16").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("C23").Value = 2
$ws.Range("E23").Value = 100
$ws.Range("F23").Value = 3
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = 50
$ws.Range("I23").Value = 19
$ws.Range("J23").Value = 23
$ws.Range("K23").Value = -17.391304347826
$ws.Range("L23").Value = 58.333333333333
$ws.Range("M23").Value = 72.727272727272
$ws.Range("C24").Value = 19
$ws.Range("D24").Value = 32
$ws.Range("E24").Value = -40.625
$ws.Range("G24").Value = 133
$ws.Range("H24").Value = -52.631578947368
$ws.Range("I24").Value = 725
$ws.Range("J24").Value = 901
$ws.Range("K24").Value = -19.533851276359
$ws.Range("L24").Value = 16.559485530546
$ws.Range("M24").Value = 2.401129943502
$ws.Range("C25").Value = 16
$ws.Range("D25").Value = 16
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 86
$ws.Range("G25").Value = 76
$ws.Range("H25").Value = 13.157894736842
$ws.Range("I25").Value = 656
$ws.Range("J25").Value = 657
$ws.Range("K25").Value = -0.152207001522
$ws.Range("L25").Value = 25.190839694656
$ws.Range("M25").Value = -8.888888888888
$ws.Range("C16").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("C26").Value = 1
$ws.Range("C14").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E26").PasteSpecial(-4163)
$ws.Range("E14").Copy()
$ws.Range("E26").PasteSpecial(-4122)
$ws.Range("G26").Value = 4
$ws.Range("H26").Value = -25
$ws.Range("I26").Value = 37
$ws.Range("K26").Value = -15.90909090909
$ws.Range("L26").Value = -9.756097560975
$ws.Range("C16").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("C27").Value = 2
$ws.Range("C16").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("D27").Value = 1
$ws.Range("K14").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("E27").Value = 100
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = -33.333333333333
$ws.Range("I27").Value = 53
$ws.Range("J27").Value = 58
$ws.Range("K27").Value = -8.620689655172
$ws.Range("L27").Value = 39.473684210526
$ws.Range("C16").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("D28").Value = 2
$ws.Range("K14").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("E28").Value = -50
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = -25
$ws.Range("I28").Value = 28
$ws.Range("J28").Value = 44
$ws.Range("K28").Value = -36.363636363636
$ws.Range("L28").Value = -37.777777777777
$ws.Range("M28").Value = -9.677419354838
$ws.Range("N28").Value = -69.565217391304
$ws.Range("C16").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("D29").Value = 2
$ws.Range("K14").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("E29").Value = -50
$ws.Range("G29").Value = 4
$ws.Range("H29").Value = -25
$ws.Range("I29").Value = 21
$ws.Range("J29").Value = 38
$ws.Range("K29").Value = -44.736842105263
$ws.Range("L29").Value = -48.780487804878
$ws.Range("M29").Value = -19.230769230769
$ws.Range("N29").Value = -73.75

$excel.CutCopyMode = $false
